$d = $word.ActiveDocument

# Locate and remove the existing _GoBack bookmark; remember where it was.
$bm = $d.Bookmarks("_GoBack")
$insertPos = $bm.Start
$bm.Delete()

# Insert the new content (13 paragraphs) as raw OOXML right where the bookmark was.
# The bookmark itself is re-created inside the final inserted paragraph, in the
# same place (immediately after the last run, before that paragraph's mark).
$insertRange = $d.Range($insertPos, $insertPos)
$newContentXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p><w:r><w:t>Service possibilities:</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Payment “database” -&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Integromat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> writes to an excel sheet, to keep record of all transactions</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Receipt -&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Integromat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> sends an email to the guest or a notification to the application regarding his order. </w:t></w:r><w:r><w:t>Non onboarded users</w:t></w:r><w:r><w:t xml:space="preserve"> can choose if they want a print-out version or a digital one (email)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Decision Support System -&gt; mongoose aggregation on database, to return desired rows according to input. </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Incentive for feedback -&gt; offer coupons for feedback, like 5% on the next order, or after 20 recommendations receive one coupon with 50% off or a free meal.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t>Rating helps overall quality of the recommendation system -&gt; desired</w:t></w:r><w:r><w:t>!</w:t></w:r></w:p><w:p/><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Possibility (for preorders) to use any voice recognition service </w:t></w:r><w:r><w:t>(</w:t></w:r><w:r><w:t xml:space="preserve">e.g. </w:t></w:r><w:r><w:t xml:space="preserve">google assistant) </w:t></w:r><w:r><w:t xml:space="preserve">to place the order, as well as </w:t></w:r><w:r><w:t>to update the restaurant on an upcoming delay.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
[void]$insertRange.InsertXML($newContentXml)

# The bookmark now sits right after the text of the newly-inserted final
# paragraph (i.e. right before that paragraph's own mark) -- use it to find
# that position without depending on a hard-coded paragraph index.
$bmAfter = $d.Bookmarks("_GoBack")
$lastContentParaEnd = $bmAfter.Start

# The six trailing empty paragraphs that used to follow the "Preorder..."
# paragraph must be removed so the new last paragraph (with the bookmark)
# becomes the final paragraph of the document body. The very last paragraph
# mark of the body cannot itself be deleted, so repeatedly deleting the
# character right after our new content collapses all of the old empty
# paragraphs into that unremovable final mark.
for ($i = 0; $i -lt 6; $i++) {
    $killRange = $d.Range($lastContentParaEnd, $lastContentParaEnd + 1)
    $killRange.Delete()
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
